# Sprout test fixes: rename a handful of template-function names inside the
# doc-gen-test fixture text, tighten up a field reference, and tidy the
# spacing on the {{range}} block paragraphs.
#
# Note: Find/Replace in this host applies smart-quote autocorrect to any
# literal `"` that appears in the *replacement* string, so every
# replacement below is written to avoid touching quote characters (the
# existing straight quotes in the document are left exactly where they
# are by only matching/replacing the text around them).

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# 1. "-7d" -> "-168h" in the dateModify example (quotes untouched).
Replace-Text '-7d' '-168h'

# 2-4. upper/lower/title -> toUpper/toLower/toTitleCase in the String
#      Functions section (these are the first occurrences; the "Customer:"
#      upper further down is handled separately below).
Replace-Text '{{.record.name | upper}}                ' '{{.record.name | toUpper}}                '
Replace-Text '{{.record.name | lower}}                ' '{{.record.name | toLower}}                '
Replace-Text '{{.record.name | title}}                ' '{{.record.name | toTitleCase}}                '

# 5. abbrev 50 -> ellipsis 50
Replace-Text '{{.record.name | abbrev 50}}' '{{.record.name | ellipsis 50}}'

# 6. float64 -> toFloat64
Replace-Text '{{.record.total_unit_price | float64}}' '{{.record.total_unit_price | toFloat64}}'

# 7. int -> toInt (Number Functions single-filter example), and the
#    trailing "Convert to int" label also becomes "Convert to toInt".
Replace-Text '{{.record.total_units | int}}' '{{.record.total_units | toInt}}'
Replace-Text '             → Convert to int' '             → Convert to toInt'

# 8. int | round 2 -> toInt | round 2
Replace-Text '{{.record.total_units | int | round 2}}' '{{.record.total_units | toInt | round 2}}'

# 9. Customer: {{.record.name | upper}} -> toUpper
Replace-Text 'Customer: {{.record.name | upper}}' 'Customer: {{.record.name | toUpper}}'

# 10. Total Value field switches from total_units to total_unit_price (an
#     extra space is kept ahead of the pipe, matching the source fixture).
Replace-Text 'total_units | printf' 'total_unit_price  | printf'

# 11-13. The {{range}} / item / {{end}} paragraphs lose their leading
#        indentation and pick up tighter spacing (0pt after, 1.5 line).
Replace-Text '  {{range .record.units}}' '{{range .record.units}}'
Replace-Text '  - {{.name}}: ${' '{{.name}}: ${'
Replace-Text '  {{end}}' '{{end}}'

function Set-TightSpacing($para) {
    $para.Format.SpaceAfter = 0
    $para.Format.LineSpacingRule = 5
    $para.Format.LineSpacing = 18
}

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -eq "{{range .record.units}}`r") {
        Set-TightSpacing $p
    } elseif ($t -eq "{{.name}}: `${{.unit_price}} x {{.additional_price}} = `${{mul .unit_price .additional_price}}`r") {
        Set-TightSpacing $p
    } elseif ($t -eq "{{end}}`r") {
        Set-TightSpacing $p
    }
}
